# "adding averages and more checks"
#
# - Training Dashboard: "PERIOD TO EXPIRE" (col H) drops by 8 days and
#   "LAST UPDATE" (col I) moves forward from 08-Sep-2025 to 16-Sep-2025
#   for every data row (3-8).
# - Exam Dashboard: widen the COMMENTS column and replace the terse "OK"
#   remarks with a clearer "date is valid" message.
# - Styling: the title drops its own oversized (14pt) font, and both the
#   title and the header-row banner text become bold white (matching the
#   dark blue banner fill).

$wb = $excel.ActiveWorkbook

$training = $wb.Worksheets.Item("Training Dashboard")
$exam = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: refresh "PERIOD TO EXPIRE" / "LAST UPDATE" ---
$periodUpdates = @{
    3 = 638
    4 = 364
    5 = 587
    6 = 423
    7 = 587
    8 = 181
}
foreach ($row in $periodUpdates.Keys) {
    $training.Cells.Item($row, 8).Value = $periodUpdates[$row]
    $lastUpdateCell = $training.Cells.Item($row, 9)
    $lastUpdateCell.NumberFormat = "@"
    $lastUpdateCell.Value = "16-Sep-2025"
}

# --- Exam Dashboard: wider COMMENTS column + friendlier remarks ---
$exam.Range("E1").EntireColumn.ColumnWidth = 14.16666666666667
$exam.Range("E3").Value = "date is valid"
$exam.Range("E4").Value = "date is valid"

# --- Styling tweaks shared by both dashboards ---
# Both the title and the header-row banner end up sharing one bold white
# font (no more separate oversized 14pt title font / plain bold header font).
foreach ($ws in @($training, $exam)) {
    $titleCell = $ws.Range("A1")
    $titleCell.Font.Size = 11
    $titleCell.Font.Bold = $true
    $titleCell.Font.Color = 16777215

    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $ws.UsedRange.Columns.Count))
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215
}
